$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.72
$ws.Range("D4").Value = 0.73
$ws.Range("D8").Value = 0.9800000000000001
$ws.Range("E8").Value = 0.8311996319189094
$ws.Range("G8").Value = 0.9151162790697674
$ws.Range("H8").Value = 0.9195402298850575
$ws.Range("D9").Value = 0.9299999999999999
$ws.Range("D11").Value = 0.9299999999999999
$ws.Range("D12").Value = 0.71
$ws.Range("D13").Value = 0.8699999999999999
$ws.Range("D15").Value = 0.9810000000000001
$ws.Range("D16").Value = 0.9810000000000001
$ws.Range("D17").Value = 0.991
$ws.Range("E17").Value = 0.9047067186771327
$ws.Range("G17").Value = 0.9526162790697674
$ws.Range("H17").Value = 0.9523809523809524
$ws.Range("D19").Value = 0.982
$ws.Range("D20").Value = 0.95
$ws.Range("D21").Value = 0.8299999999999998
$ws.Range("D22").Value = 0.9199999999999999
$ws.Range("E22").Value = 0.9074162525099211
$ws.Range("G22").Value = 0.95
$ws.Range("H22").Value = 0.9555555555555556
$ws.Range("D24").Value = 0.994
$ws.Range("D25").Value = 0.994
$ws.Range("D26").Value = 0.994
$ws.Range("D28").Value = 0.993
$ws.Range("E28").Value = 0.9278236787666257
$ws.Range("G28").Value = 0.9633720930232558
$ws.Range("H28").Value = 0.9655172413793104
$ws.Range("D29").Value = 0.9970000000000001
$ws.Range("D36").Value = 0.82
$ws.Range("D37").Value = 0.82
$ws.Range("E37").Value = 0.7952302631578947
$ws.Range("G37").Value = 0.8976151315789473
$ws.Range("H37").Value = 0.953125
$ws.Range("D41").Value = 0.85
$ws.Range("D43").Value = 0.8400000000000001
$ws.Range("E43").Value = 0.931463152500253
$ws.Range("F43").Value = 0.9759036144578314
$ws.Range("G43").Value = 0.9473684210526316
$ws.Range("H43").Value = 0.9846153846153847
$ws.Range("D45").Value = 0.72
$ws.Range("E45").Value = 0.6472299111123082
$ws.Range("G45").Value = 0.8108552631578947
$ws.Range("H45").Value = 0.923076923076923
$ws.Range("E46").Value = 0.3298551901380059
$ws.Range("F46").Value = 0.7349397590361446
$ws.Range("G46").Value = 0.6800986842105263
$ws.Range("H46").Value = 0.8196721311475409
$ws.Range("E47").Value = 0.5638878933294721
$ws.Range("F47").Value = 0.8554216867469879
$ws.Range("G47").Value = 0.7582236842105263
$ws.Range("H47").Value = 0.9090909090909091
$ws.Range("E51").Value = 0.7514898301685916
$ws.Range("F51").Value = 0.9156626506024096
$ws.Range("G51").Value = 0.852796052631579
$ws.Range("H51").Value = 0.9465648854961832
$ws.Range("D52").Value = 0.9
$ws.Range("D55").Value = 0.8699999999999999
$ws.Range("E55").Value = 0.699828204523419
$ws.Range("F55").Value = 0.9036144578313253
$ws.Range("G55").Value = 0.8850931677018633
$ws.Range("H55").Value = 0.75
$ws.Range("D57").Value = 0.7899999999999999
$ws.Range("E57").Value = 0.699828204523419
$ws.Range("F57").Value = 0.9036144578313253
$ws.Range("G57").Value = 0.8850931677018633
$ws.Range("H57").Value = 0.75
$ws.Range("D58").Value = 0.97
$ws.Range("E58").Value = 0.8681907322380075
$ws.Range("G58").Value = 0.9213250517598344
$ws.Range("H58").Value = 0.888888888888889
$ws.Range("D59").Value = 0.986
$ws.Range("D62").Value = 0.9869999999999999
$ws.Range("D63").Value = 0.9099999999999999
$ws.Range("D64").Value = 0.86
$ws.Range("D65").Value = 0.982
$ws.Range("D66").Value = 0.9869999999999999
$ws.Range("E66").Value = 0.8281573498964804
$ws.Range("G66").Value = 0.9140786749482401
$ws.Range("H66").Value = 0.8571428571428571
$ws.Range("D68").Value = 0.9890000000000001
$ws.Range("D70").Value = 0.99
$ws.Range("D72").Value = 0.9949999999999999
$ws.Range("D74").Value = 0.9869999999999999
$ws.Range("D75").Value = 0.992
$ws.Range("D76").Value = 0.9800000000000001
$ws.Range("D77").Value = 0.95
$ws.Range("E77").Value = 0.7918335465004855
$ws.Range("G77").Value = 0.9068322981366459
$ws.Range("H77").Value = 0.8275862068965518
$ws.Range("D78").Value = 0.7
$ws.Range("E78").Value = 0.551946004289372
$ws.Range("F78").Value = 0.8433734939759037
$ws.Range("G78").Value = 0.8203933747412009
$ws.Range("H78").Value = 0.6285714285714286
$ws.Range("D80").Value = 0.9099999999999999
$ws.Range("E80").Value = 0.6840677729207231
$ws.Range("F80").Value = 0.9397590361445783
$ws.Range("G80").Value = 0.967948717948718
$ws.Range("H80").Value = 0.6666666666666666
$ws.Range("D85").Value = 0.9299999999999999
$ws.Range("D86").Value = 0.992
$ws.Range("E86").Value = 0.7871794871794872
$ws.Range("F86").Value = 0.9759036144578314
$ws.Range("G86").Value = 0.8935897435897436
$ws.Range("H86").Value = 0.8000000000000002
$ws.Range("D87").Value = 0.982
$ws.Range("D89").Value = 0.9890000000000001
$ws.Range("D90").Value = 0.99
$ws.Range("D94").Value = 0.9869999999999999
$ws.Range("D97").Value = 0.99
$ws.Range("D99").Value = 0.9890000000000001
$ws.Range("D100").Value = 0.982
$ws.Range("D110").Value = 0.85
$ws.Range("E110").Value = 0.4481697548529493
$ws.Range("F110").Value = 0.7108433734939759
$ws.Range("G110").Value = 0.7185314685314685
$ws.Range("H110").Value = 0.6842105263157896
$ws.Range("D117").Value = 0.7
$ws.Range("E117").Value = 0.5878821097951015
$ws.Range("G117").Value = 0.7896270396270396
$ws.Range("H117").Value = 0.7692307692307693
$ws.Range("D120").Value = 0.9949999999999999
$ws.Range("E120").Value = 0.6036949306884013
$ws.Range("F120").Value = 0.7831325301204819
$ws.Range("G120").Value = 0.7706876456876457
$ws.Range("H120").Value = 0.826923076923077
$ws.Range("D125").Value = 0.9399999999999999
$ws.Range("E125").Value = 0.3249165681458976
$ws.Range("F125").Value = 0.8915662650602409
$ws.Range("G125").Value = 0.7833333333333333
$ws.Range("H125").Value = 0.3076923076923077
$ws.Range("D126").Value = 0.93
$ws.Range("D127").Value = 0.96
$ws.Range("E127").Value = 0.2578155966122687
$ws.Range("F127").Value = 0.9397590361445783
$ws.Range("G127").Value = 0.6479166666666667
$ws.Range("H127").Value = 0.2857142857142858
$ws.Range("D128").Value = 0.983
$ws.Range("D129").Value = 0.997
$ws.Range("D130").Value = 0.988
$ws.Range("E130").Value = 0.4057948819484283
$ws.Range("F130").Value = 0.927710843373494
$ws.Range("G130").Value = 0.8020833333333333
$ws.Range("H130").Value = 0.4
$ws.Range("D131").Value = 0.89
$ws.Range("D132").Value = 0.96
$ws.Range("E132").Value = 0.4444624815964758
$ws.Range("F132").Value = 0.9397590361445783
$ws.Range("G132").Value = 0.8083333333333333
$ws.Range("H132").Value = 0.4444444444444444
$ws.Range("D134").Value = 0.9399999999999999
$ws.Range("E134").Value = 0.2578155966122687
$ws.Range("F134").Value = 0.9397590361445783
$ws.Range("G134").Value = 0.6479166666666667
$ws.Range("H134").Value = 0.2857142857142858
$ws.Range("D135").Value = 0.9899999999999999
$ws.Range("E135").Value = 0.4444624815964758
$ws.Range("F135").Value = 0.9397590361445783
$ws.Range("G135").Value = 0.8083333333333333
$ws.Range("H135").Value = 0.4444444444444444
$ws.Range("D137").Value = 0.96
$ws.Range("E137").Value = 0.2717549127444297
$ws.Range("F137").Value = 0.8554216867469879
$ws.Range("G137").Value = 0.7645833333333334
$ws.Range("H137").Value = 0.25
$ws.Range("D138").Value = 0.8799999999999999
$ws.Range("E138").Value = 0.1685491513628021
$ws.Range("F138").Value = 0.7349397590361446
$ws.Range("G138").Value = 0.7020833333333334
$ws.Range("H138").Value = 0.1538461538461538
$ws.Range("D139").Value = 0.989
$ws.Range("E139").Value = 0.2575298062442672
$ws.Range("F139").Value = 0.8433734939759037
$ws.Range("G139").Value = 0.7583333333333333
$ws.Range("H139").Value = 0.2352941176470588
$ws.Range("E140").Value = 0.1348800974667526
$ws.Range("F140").Value = 0.6746987951807228
$ws.Range("G140").Value = 0.6708333333333334
$ws.Range("H140").Value = 0.1290322580645161
$ws.Range("D141").Value = 0.9899999999999999
$ws.Range("E141").Value = 0.4057948819484283
$ws.Range("F141").Value = 0.927710843373494
$ws.Range("G141").Value = 0.8020833333333333
$ws.Range("H141").Value = 0.4
